$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text values that look like plain numbers need to be forced to stay text
# (matching the workbooks convention of storing prices as inline strings),
# otherwise Excel auto-converts them to numeric values on assignment.

$ws.Range('D2').Value = '51.008.31'
$ws.Range('D3').Value = '2.944.92'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '375.58'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.31%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '101.14'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.30%  '
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.38'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.28%  '
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '3.404.59'
$ws.Range('E13').Value = '  -1.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.15'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '11.28'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +51.39%  '
$ws.Range('D17').Value = '2.946.77'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.998'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').Value = '50.966.06'
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('E20').Value = '  -6.02%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.49'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.26%  '
$ws.Range('D22').Value = '0.0₃0956'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '265.50'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.87'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.15'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +7.33%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.14'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.13%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.56'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '25.72'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.164'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.83%  '
$ws.Range('E32').Value = '  -4.96%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '10.01'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '50.94'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '33.42'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.15%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0443'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.16'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.77%  '
$ws.Range('E40').Value = '  -0.59%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '16.51'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.28%  '
$ws.Range('E42').Value = '  -1.90%  '
$ws.Range('E43').Value = '  -5.03%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '120.52'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.26'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.71%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.42'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.30%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.272'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.45%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.33'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').Value = '1.992.28'
$ws.Range('E50').Value = '  -2.22%  '
$ws.Range('E51').Value = '  -1.75%  '
